# Natmi following Dr Hou advice:
# Re-run of the NATMI ligand-receptor summary for Rspo2-Lgr4 (Sending cluster FAPs) now also
# includes "sCs" as a sending/target cluster and "M2" as a target cluster, expanding the
# result table from 3 data rows (rows 2-4) to 8 data rows (rows 2-9) and refreshing the
# computed expression/specificity statistics for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.597878666666666
$ws.Range("H2").Value = 4.793635999999999
$ws.Range("I2").Value = 0.98224549682877
$ws.Range("J2").Value = 0.9822454968287699
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.369012666666667
$ws.Range("N2").Value = 4.107038
$ws.Range("O2").Value = 0.06074872832285554
$ws.Range("P2").Value = 0.06074872832285555
$ws.Range("Q2").Value = 2.187516134463111
$ws.Range("R2").Value = 19.687645210168
$ws.Range("S2").Value = 0.05967016483319922
$ws.Range("T2").Value = 0.05967016483319922

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.597878666666666
$ws.Range("H3").Value = 4.793635999999999
$ws.Range("I3").Value = 0.98224549682877
$ws.Range("J3").Value = 0.9822454968287699
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.382866
$ws.Range("N3").Value = 43.148598
$ws.Range("O3").Value = 0.6382269794957115
$ws.Range("P3").Value = 0.6382269794957116
$ws.Range("Q3").Value = 22.98207474692533
$ws.Range("R3").Value = 206.838672722328
$ws.Range("S3").Value = 0.6268955765642904
$ws.Range("T3").Value = 0.6268955765642904

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.597878666666666
$ws.Range("H4").Value = 4.793635999999999
$ws.Range("I4").Value = 0.98224549682877
$ws.Range("J4").Value = 0.9822454968287699
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05862
$ws.Range("N4").Value = 0.17586
$ws.Range("O4").Value = 0.002601210741867345
$ws.Range("P4").Value = 0.002601210741867345
$ws.Range("Q4").Value = 0.09366764743999999
$ws.Range("R4").Value = 0.8430088269599998
$ws.Range("S4").Value = 0.002555027537501824
$ws.Range("T4").Value = 0.002555027537501824

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.597878666666666
$ws.Range("H5").Value = 4.793635999999999
$ws.Range("I5").Value = 0.98224549682877
$ws.Range("J5").Value = 0.9822454968287699
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.725161
$ws.Range("N5").Value = 20.175483
$ws.Range("O5").Value = 0.2984230814395656
$ws.Range("P5").Value = 0.2984230814395656
$ws.Range("Q5").Value = 10.74599129179867
$ws.Range("R5").Value = 96.71392162618798
$ws.Range("S5").Value = 0.2931247278937786
$ws.Range("T5").Value = 0.2931247278937786

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rspo2"
$ws.Range("C6").Value = "Lgr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02888233333333333
$ws.Range("H6").Value = 0.086647
$ws.Range("I6").Value = 0.01775450317123003
$ws.Range("J6").Value = 0.01775450317123003
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.369012666666667
$ws.Range("N6").Value = 4.107038
$ws.Range("O6").Value = 0.06074872832285554
$ws.Range("P6").Value = 0.06074872832285555
$ws.Range("Q6").Value = 0.03954028017622222
$ws.Range("R6").Value = 0.355862521586
$ws.Range("S6").Value = 0.00107856348965633
$ws.Range("T6").Value = 0.001078563489656331

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rspo2"
$ws.Range("C7").Value = "Lgr4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02888233333333333
$ws.Range("H7").Value = 0.086647
$ws.Range("I7").Value = 0.01775450317123003
$ws.Range("J7").Value = 0.01775450317123003
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.382866
$ws.Range("N7").Value = 43.148598
$ws.Range("O7").Value = 0.6382269794957115
$ws.Range("P7").Value = 0.6382269794957116
$ws.Range("Q7").Value = 0.4154107301006666
$ws.Range("R7").Value = 3.738696570906
$ws.Range("S7").Value = 0.01133140293142117
$ws.Range("T7").Value = 0.01133140293142118

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rspo2"
$ws.Range("C8").Value = "Lgr4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.02888233333333333
$ws.Range("H8").Value = 0.086647
$ws.Range("I8").Value = 0.01775450317123003
$ws.Range("J8").Value = 0.01775450317123003
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05862
$ws.Range("N8").Value = 0.17586
$ws.Range("O8").Value = 0.002601210741867345
$ws.Range("P8").Value = 0.002601210741867345
$ws.Range("Q8").Value = 0.00169308238
$ws.Range("R8").Value = 0.01523774142
$ws.Range("S8").Value = 0.00004618320436552139
$ws.Range("T8").Value = 0.0000461832043655214

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rspo2"
$ws.Range("C9").Value = "Lgr4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.02888233333333333
$ws.Range("H9").Value = 0.086647
$ws.Range("I9").Value = 0.01775450317123003
$ws.Range("J9").Value = 0.01775450317123003
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.725161
$ws.Range("N9").Value = 20.175483
$ws.Range("O9").Value = 0.2984230814395656
$ws.Range("P9").Value = 0.2984230814395656
$ws.Range("Q9").Value = 0.1942383417223333
$ws.Range("R9").Value = 1.748145075501
$ws.Range("S9").Value = 0.005298353545787005
$ws.Range("T9").Value = 0.005298353545787005

